# Optuna Attempt (go back with original)
# Updates Seasonality Index / Inventory Coverage / MyForecast / Reorder Urgency values
# on the "Forecast Comparison" sheet, and the two derived totals on "Summary".

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------

# Seasonality Index (column L)
$wsForecast.Range("L2").Value  = 0.91
$wsForecast.Range("L3").Value  = 1.09
$wsForecast.Range("L4").Value  = 0.95
$wsForecast.Range("L5").Value  = 0.82
$wsForecast.Range("L6").Value  = 0.97
$wsForecast.Range("L7").Value  = 1.02
$wsForecast.Range("L8").Value  = 0.99
$wsForecast.Range("L9").Value  = 0.84
$wsForecast.Range("L10").Value = 1.08
$wsForecast.Range("L11").Value = 1.04
$wsForecast.Range("L12").Value = 1.18
$wsForecast.Range("L13").Value = 0.89
$wsForecast.Range("L14").Value = 1.14
$wsForecast.Range("L15").Value = 1
$wsForecast.Range("L16").Value = 0.93
$wsForecast.Range("L17").Value = 1.15

# MyForecast (column D)
$wsForecast.Range("D10").Value = 13
$wsForecast.Range("D12").Value = 13
$wsForecast.Range("D13").Value = 13
$wsForecast.Range("D14").Value = 12
$wsForecast.Range("D15").Value = 11
$wsForecast.Range("D16").Value = 11

# Inventory Coverage (column H)
$wsForecast.Range("H10").Value = 1.8
$wsForecast.Range("H11").Value = 0.89
$wsForecast.Range("H12").Value = 0

# Reorder Urgency (column J)
$wsForecast.Range("J11").Value = "Urgent"

# --- Summary sheet --------------------------------------------------------------
# B9 and B14 are stored as text (numeric-looking strings), so use a leading
# apostrophe (quote prefix) to keep them as text cells rather than numbers.

$wsSummary.Range("B9").Formula  = "'203"
$wsSummary.Range("B14").Formula = "'10"
